$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.186.88"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "1.831.50"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("E4").Value = "  +0.42%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.63"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3666"
$ws.Range("E8").Value = "  -0.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07410"
$ws.Range("E9").Value = "  +0.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8814"
$ws.Range("E10").Value = "  +0.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.37"
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").Value = "1.889.24"
$ws.Range("E12").Value = "  +2.98%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07675"
$ws.Range("E13").Value = "  +5.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.392"
$ws.Range("E14").Value = "  -1.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.31"
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.545"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.009"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008745"
$ws.Range("E18").Value = "  -0.29%  "
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("D20").Value = "27.584.39"
$ws.Range("E20").Value = "  +1.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.65"
$ws.Range("E21").Value = "  -0.73%  "
$ws.Range("E22").Value = "  -1.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.62"
$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("D24").Value = "2.090.05"
$ws.Range("E24").Value = "  +1.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.885"
$ws.Range("E25").Value = "  -0.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.18"
$ws.Range("E26").Value = "  -0.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.53"
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("E28").Value = "  -1.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.192"
$ws.Range("E29").Value = "  -1.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.71"
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08938"
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7462"
$ws.Range("E32").Value = "  -1.53%  "
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.962"
$ws.Range("E34").Value = "  +1.25%  "
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("E36").Value = "  +0.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.538"
$ws.Range("E37").Value = "  +5.07%  "
$ws.Range("E38").Value = "  -0.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05303"
$ws.Range("E39").Value = "  -0.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01939"
$ws.Range("E40").Value = "  -0.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.309"
$ws.Range("E41").Value = "  +0.55%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.943"
$ws.Range("E42").Value = "  -1.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5278"
$ws.Range("E43").Value = "  -1.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1644"
$ws.Range("E44").Value = "  -1.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.397"
$ws.Range("E45").Value = "  -1.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4910"
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.37"
$ws.Range("E47").Value = "  -1.55%  "
$ws.Range("E48").Value = "  +0.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.60"
$ws.Range("E49").Value = "  +0.94%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.656"
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06281"
$ws.Range("E51").Value = "  -0.45%  "
